$d = $word.ActiveDocument

# The document has three "<id>...</id>" fields, each split across three
# runs: "<id>" (Courier New / 7f6000), "p083v_aN" (plain black), "</id>"
# (Courier New / 7f6000). The edit collapses each triple into a single
# run reading "<id>p083v_N</id>" (using the formatting of the first run),
# for N = 1, 2, 3.

$d.Content.Find.Execute("<id>p083v_a1</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p083v_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p083v_a2</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p083v_2</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p083v_a3</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p083v_3</id>", 2) | Out-Null
